# Replace the hard-coded "MEXICANA" nationality in the two affidavit
# paragraphs with the ${NACIONALIDAD_ARRENDADOR} / ${NACIONALIDAD_ARRENDATARIO}
# placeholders, splitting the original single run into three runs
# (prefix text, placeholder, trailing period) while keeping the run
# formatting (rFonts/color/kern/sz/szCs/lang) identical on all three.

$d = $word.ActiveDocument

function Replace-Nacionalidad($anchorText, $placeholder) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $anchorText"
        return
    }

    # Replace just the word "MEXICANA" with the placeholder text first; since
    # the formatting is unchanged this merges back into the surrounding run.
    $rng.Text = $placeholder

    # Touching a character-level formatting property on the (now resized)
    # found range forces the engine to split the paragraph's run into three
    # runs: the unchanged leading text, the placeholder itself, and the
    # trailing "." -- all three keep identical rPr because we restore Bold
    # to its original (unset/false) value.
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

Replace-Nacionalidad "MEXICANA" "`${NACIONALIDAD_ARRENDADOR}"
Replace-Nacionalidad "MEXICANA" "`${NACIONALIDAD_ARRENDATARIO}"
